$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as plain text (preserving formatting like trailing
# zeros / thousands separators that would otherwise be lost if Excel
# auto-converted the string to a number).
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Set-Row($row, $d, $e) {
    if ($null -ne $d) {
        Set-TextCell $row 4 $d
    }
    if ($null -ne $e) {
        Set-TextCell $row 5 $e
    }
}

Set-Row 2  "60.702.66"  "  -1.89%  "
Set-Row 3  "3.383.64"   "  -2.29%  "
Set-Row 4  $null         "  +0.07%  "
Set-Row 5  "569.80"     "  -2.36%  "
Set-Row 6  "141.56"     "  -4.08%  "
Set-Row 7  $null         "  +0.08%  "
Set-Row 8  "3.383.29"   "  -2.36%  "
Set-Row 9  $null         "  -0.29%  "
Set-Row 10 "7.52"       "  -1.95%  "
Set-Row 11 $null         "  -2.40%  "
Set-Row 12 "0.401"      "  +2.09%  "
Set-Row 13 "3.965.38"   "  -2.13%  "
Set-Row 14 "28.42"      "  +1.36%  "
Set-Row 15 $null         "  +1.53%  "
Set-Row 16 $null         "  -2.38%  "
Set-Row 17 "3.385.79"   "  -2.23%  "
Set-Row 18 "60.814.28"  "  -1.84%  "
Set-Row 19 $null         "  +0.02%  "
Set-Row 20 "14.04"      "  -2.48%  "
Set-Row 21 "9.01"       "  -6.05%  "
Set-Row 22 "384.68"     "  -1.38%  "
Set-Row 23 "0.562"      "  -0.93%  "
Set-Row 24 "73.65"      "  -0.19%  "
Set-Row 25 "1.00"       "  +0.20%  "
Set-Row 26 $null         "  -5.92%  "
Set-Row 27 "3.522.76"   "  -2.14%  "
Set-Row 28 $null         "  -2.53%  "
Set-Row 29 "0.998"      "  +0.22%  "
Set-Row 30 "7.42"       "  -4.13%  "

# Rows 31-33 got reordered (InternetComputer, PancakeSwap, Fetch.AI)
# -> (Fetch.AI, InternetComputer, PancakeSwap), with updated values.
Set-TextCell 31 2 "Fetch.AI"
Set-TextCell 31 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell 31 4 "1.44"
Set-TextCell 31 5 "  -3.15%  "

Set-TextCell 32 2 "InternetComputer(DFINITY)"
Set-TextCell 32 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 32 4 "8.00"
Set-TextCell 32 5 "  -2.98%  "

Set-TextCell 33 2 "PancakeSwap"
Set-TextCell 33 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell 33 4 "2.15"
Set-TextCell 33 5 "  -4.34%  "

Set-Row 34 $null         "  -0.03%  "
Set-Row 35 $null         "  -2.38%  "
Set-Row 36 "6.99"       "  -0.51%  "
Set-Row 37 "166.45"     "  -0.31%  "

# Rows 38-39 got reordered (RenzoRestakedETH, NEARProtocol)
# -> (NEARProtocol, RenzoRestakedETH), with updated values.
Set-TextCell 38 2 "NEARProtocol"
Set-TextCell 38 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 38 4 "5.01"
Set-TextCell 38 5 "  -3.17%  "

Set-TextCell 39 2 "RenzoRestakedETH"
Set-TextCell 39 3 "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextCell 39 4 "3.416.05"
Set-TextCell 39 5 "  -2.12%  "

Set-Row 40 $null         "  -5.24%  "

# Rows 41-42 got reordered (EnergySwap, Hedera) -> (Hedera, EnergySwap),
# with updated values.
Set-TextCell 41 2 "Hedera"
Set-TextCell 41 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 41 4 "0.0779"
Set-TextCell 41 5 "  -0.91%  "

Set-TextCell 42 2 "EnergySwap"
Set-TextCell 42 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 42 4 "27.68"
Set-TextCell 42 5 "  +0.92%  "

Set-Row 43 $null         "  -3.33%  "
Set-Row 44 $null         "  +0.12%  "
Set-Row 45 "4.44"       "  -1.71%  "
Set-Row 46 "41.68"      "  -2.05%  "
Set-Row 47 $null         "  -3.15%  "
Set-Row 48 "2.532.56"   "  -1.62%  "
Set-Row 49 "1.12"       "  -4.59%  "
Set-Row 50 "23.52"      "  +1.54%  "
Set-Row 51 "6.88"       "  -1.09%  "
